$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: split "Terms Typically Offered" semantics into new columns
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Row 2
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"

# Row 3
$ws.Range("C3").Value = "CRP 201, CRP 211 or consent of instructor."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "W"

# Row 4
$ws.Range("C4").Value = "CRP 202 and CRP 204."
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "SP"

# Row 5
$ws.Range("C5").Value = "CRP 201."
$ws.Range("D5").Value = "CRP 202."
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "W "

# Row 6
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "W"

# Row 7
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "F, W, SP"

# Row 8
$ws.Range("C8").Value = "CRP 212."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F"

# Row 9
$ws.Range("C9").Value = "CRP 212."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "SP"

# Row 10
$ws.Range("C10").Value = "Completion of GE Area D1."
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "ES 112."
$ws.Range("G10").Value = "SP "

# Row 11
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "F"

# Row 12
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "TBD"

# Row 13
$ws.Range("C13").Value = "CRP 212."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "W"

# Row 14
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "W"

# Row 15
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "SP"

# Row 16
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "W"

# Row 17
$ws.Range("C17").Value = "CRP 212."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "F"

# Row 18
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "W, SP"

# Row 19
$ws.Range("C19").Value = "CRP 203, CRP 213, or CRP 214."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "SP"

# Row 20
$ws.Range("C20").Value = "CRP 336."
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "W"

# Row 21
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "TBD"

# Row 22
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "F, W, SP"

# Row 23
$ws.Range("C23").Value = "CRP/NR 351."
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "TBD"

# Row 24
$ws.Range("C24").Value = "Completion of GE Area A1 with a grade of C- or better and CRP 341; or graduate standing."
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "TBD"

# Row 25
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "SP"

# Row 26
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "W"

# Row 27
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "F, SP"

# Row 28
$ws.Range("C28").Value = "CRP 336, CRP 341 or consent of instructor."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F"

# Row 29
$ws.Range("C29").Value = "CRP 342, CRP 410, or consent of instructor."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "W"

# Row 30
$ws.Range("C30").Value = "CRP 212."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "SP"

# Row 31
$ws.Range("C31").Value = "CRP 212 and upper division standing."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "W"

# Row 32
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "TBD"

# Row 33
$ws.Range("C33").Value = "Junior standing."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "CRP 334."
$ws.Range("G33").Value = "F "

# Row 34
$ws.Range("C34").Value = "CRP 212."
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "F"

# Row 35
$ws.Range("C35").Value = "CRP 212, senior standing, or graduate standing."
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "F"

# Row 36
$ws.Range("C36").Value = "CRP 212 or graduate standing."
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "TBD"

# Row 37
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "TBD"

# Row 38
$ws.Range("C38").Value = "CRP 212 or graduate standing."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "CRP 336."
$ws.Range("G38").Value = "SP "

# Row 39
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "W"

# Row 40
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "TBD"

# Row 41
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "TBD"

# Row 42
$ws.Range("C42").Value = "CRP 341 or graduate standing."
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "W"

# Row 43
$ws.Range("C43").Value = "CRP 201 and CRP 202, Upper division or graduate standing."
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "TBD"

# Row 44
$ws.Range("C44").Value = "CRP 212, or graduate standing."
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "SP"

# Row 45
$ws.Range("C45").Value = "CRP 216 and junior standing, or graduate standing."
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "F"

# Row 46
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "F"

# Row 47
$ws.Range("C47").Value = "CRP 341, CRP 342."
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "F, W, SP"

# Row 48
$ws.Range("C48").Value = "CRP 410."
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "F, W, SP"

# Row 49
$ws.Range("C49").Value = "CRP 410 and senior standing."
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "SP"

# Row 50
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "TBD"

# Row 51
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "TBD"

# Row 52
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "TBD"

# Row 53
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "TBD"

# Row 54
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "F, W, SP"

# Row 55
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "F"

# Row 56
$ws.Range("D56").Value = "NA"
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "F"

# Row 57
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "TBD"

# Row 58
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "F"

# Row 59
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "W"

# Row 60
$ws.Range("D60").Value = "NA"
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "SP"

# Row 61
$ws.Range("D61").Value = "NA"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "W"

# Row 62
$ws.Range("C62").Value = "CRP 501."
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "F"

# Row 63
$ws.Range("C63").Value = "CRP 501 or consent of instructor."
$ws.Range("D63").Value = "NA"
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "NA"
$ws.Range("G63").Value = "TBD"

# Row 64
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "NA"
$ws.Range("G64").Value = "SP"

# Row 65
$ws.Range("C65").Value = "CRP 501 or graduate standing."
$ws.Range("D65").Value = "NA"
$ws.Range("E65").Value = "NA"
$ws.Range("F65").Value = "NA"
$ws.Range("G65").Value = "SP"

# Row 66
$ws.Range("D66").Value = "NA"
$ws.Range("E66").Value = "NA"
$ws.Range("F66").Value = "NA"
$ws.Range("G66").Value = "W"

# Row 67
$ws.Range("D67").Value = "NA"
$ws.Range("E67").Value = "NA"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "W"

# Row 68
$ws.Range("C68").Value = "CRP 501, CRP 525, or consent of instructor."
$ws.Range("D68").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = "NA"
$ws.Range("G68").Value = "F"

# Row 69
$ws.Range("C69").Value = "CRP 512 or consent of instructor."
$ws.Range("D69").Value = "NA"
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = "NA"
$ws.Range("G69").Value = "SP"

# Row 70
$ws.Range("C70").Value = "CRP 552."
$ws.Range("D70").Value = "NA"
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "W"

# Row 71
$ws.Range("C71").Value = "CRP 554, or consent of instructor."
$ws.Range("D71").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = "NA"
$ws.Range("G71").Value = "SP"

# Row 72
$ws.Range("D72").Value = "NA"
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "TBD"

# Row 73
$ws.Range("D73").Value = "NA"
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = "NA"
$ws.Range("G73").Value = "TBD"

# Row 74
$ws.Range("C74").Value = "CRP 513, and consent of the graduate program coordinator."
$ws.Range("D74").Value = "NA"
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = "NA"
$ws.Range("G74").Value = "F, W, SP"

# Row 75
$ws.Range("C75").Value = "CRP 513, and consent of the graduate program coordinator."
$ws.Range("D75").Value = "NA"
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = "NA"
$ws.Range("G75").Value = "F, W, SP"
